$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all edited cells so numeric-looking strings
# (e.g. "565.80", "13.30", "0.999") are preserved verbatim as text,
# matching the original inlineStr cell type instead of being
# auto-converted to numbers (which would drop trailing zeros).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.893.13'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.96%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.024.87'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.28%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '565.80'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.66'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +8.08%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.015.29'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.18%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +6.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.30'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +11.41%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000233'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.17'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +4.24%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.519.11'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.25'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +6.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.016.92'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '59.792.86'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '437.71'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +5.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.76'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.725'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +6.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.15'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.84%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.30'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.95'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.41%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.28'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +15.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.56'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.87'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +5.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.15'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.30'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +5.88%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0794'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +17.05%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +7.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.95'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +5.74%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.14'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.23'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.74%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.81'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +10.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '405.40'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +8.26%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.787.74'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.23%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.255'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +7.17%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '123.36'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.04'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.77%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.01'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +21.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.67'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.45%  '
